$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    # Force the cell to keep a literal text representation (matches the
    # original inline-string cells which must not be re-interpreted as
    # numbers by Excel's automatic type inference), then drop back to the
    # default "Normal" style so no stray number-format style sticks to the
    # cell (keeping it identical, style-wise, to its original state).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = '27.030.83'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '1.847.20'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  +0.67%  '
$ws.Range("E5").Value = '  +0.61%  '
Set-TextValue $ws "D6" '309.70'
$ws.Range("E6").Value = '  -0.14%  '
Set-TextValue $ws "D7" '0.4766'
$ws.Range("E7").Value = '  +2.02%  '
Set-TextValue $ws "D8" '0.3682'
$ws.Range("E8").Value = '  +1.63%  '
Set-TextValue $ws "D9" '0.07236'
$ws.Range("E9").Value = '  +1.28%  '
Set-TextValue $ws "D10" '0.9315'
$ws.Range("E10").Value = '  +1.89%  '
Set-TextValue $ws "D11" '19.87'
$ws.Range("E11").Value = '  +1.63%  '
Set-TextValue $ws "D12" '0.07778'
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("D13").Value = '1.859.34'
$ws.Range("E13").Value = '  -0.04%  '
Set-TextValue $ws "D14" '5.394'
$ws.Range("E14").Value = '  +2.20%  '
Set-TextValue $ws "D15" '6.478'
$ws.Range("E15").Value = '  +1.03%  '
Set-TextValue $ws "D16" '88.99'
$ws.Range("E16").Value = '  +0.73%  '
Set-TextValue $ws "D17" '1.017'
$ws.Range("E17").Value = '  +0.65%  '
Set-TextValue $ws "D18" '0.000008660'
$ws.Range("E18").Value = '  +0.83%  '
$ws.Range("E19").Value = '  +0.61%  '
$ws.Range("D20").Value = '27.035.92'
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("E21").Value = '  +1.32%  '
Set-TextValue $ws "D22" '5.052'
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("E23").Value = '  +0.08%  '
Set-TextValue $ws "D24" '1.923'
$ws.Range("E24").Value = '  -0.59%  '
Set-TextValue $ws "D25" '152.83'
$ws.Range("E25").Value = '  +0.12%  '
Set-TextValue $ws "D26" '18.29'
$ws.Range("E26").Value = '  +0.34%  '
Set-TextValue $ws "D27" '1.993'
$ws.Range("E27").Value = '  -2.38%  '
Set-TextValue $ws "D28" '114.61'
$ws.Range("E28").Value = '  +0.42%  '
Set-TextValue $ws "D29" '4.942'
$ws.Range("E29").Value = '  +0.86%  '
Set-TextValue $ws "D30" '0.08867'
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E31").Value = '  +3.61%  '
$ws.Range("E32").Value = '  +0.48%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws "D33" '0.7385'
$ws.Range("E33").Value = '  -1.17%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws "D34" '4.511'
$ws.Range("E34").Value = '  +0.98%  '
Set-TextValue $ws "D35" '2.663'
$ws.Range("E35").Value = '  -6.35%  '
Set-TextValue $ws "D36" '1.114'
$ws.Range("E36").Value = '  +2.93%  '
Set-TextValue $ws "D37" '0.01973'
$ws.Range("E37").Value = '  +1.70%  '
Set-TextValue $ws "D38" '0.05259'
$ws.Range("E38").Value = '  +1.78%  '
Set-TextValue $ws "D39" '2.968'
$ws.Range("E39").Value = '  -0.76%  '
Set-TextValue $ws "D40" '0.5283'
$ws.Range("E40").Value = '  +1.74%  '
Set-TextValue $ws "D41" '7.027'
$ws.Range("E41").Value = '  +1.73%  '
Set-TextValue $ws "D42" '0.1519'
$ws.Range("E42").Value = '  +0.46%  '
Set-TextValue $ws "D43" '8.286'
$ws.Range("E43").Value = '  +1.87%  '
Set-TextValue $ws "D44" '10.62'
$ws.Range("E44").Value = '  +0.95%  '
Set-TextValue $ws "D45" '0.4745'
$ws.Range("E45").Value = '  +0.92%  '
$ws.Range("E46").Value = '  +0.63%  '
Set-TextValue $ws "D47" '101.78'
$ws.Range("E47").Value = '  +0.95%  '
Set-TextValue $ws "D48" '1.613'
$ws.Range("E48").Value = '  +0.36%  '
Set-TextValue $ws "D49" '65.93'
$ws.Range("E49").Value = '  +1.97%  '
Set-TextValue $ws "D50" '0.06062'
Set-TextValue $ws "D51" '0.8928'
$ws.Range("E51").Value = '  +3.53%  '
